# Monday results (sorry for the delay)
# Swap the full content of row 13 and row 14 in the dataset: the lab results
# for sampleid 253 and 299 were entered under each other's rows. This
# restores them to the correct row (everything that differs between the two
# rows gets swapped; columns that already agreed between the rows are left
# untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sampleid (text column - force text so "253"/"299" aren't coerced to numbers) ---
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "253"
$ws.Range("A13").Style = "Normal"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "299"
$ws.Range("A14").Style = "Normal"

# --- date_collect ---
$ws.Range("G13").Value = 44707
$ws.Range("G14").Value = 44695

# --- sampletype ---
$ws.Range("H13").Value = "swab"
$ws.Range("H14").Value = "swab; spatula"

# --- collection ---
$ws.Range("I13").Value = "swab"
$ws.Range("I14").Value = "multiple methods"

# --- expectedsubstance ---
$ws.Range("K13").Value = "crack"
$ws.Range("K14").Value = "heroin; fentanyl"

# --- expect_opioid ---
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 1

# --- expect_fentanyl ---
$ws.Range("M13").Value = 0
$ws.Range("M14").Value = 1

# --- expect_stimulant ---
$ws.Range("O13").Value = 1
$ws.Range("O14").Value = 0

# --- expect_cocaine ---
$ws.Range("R13").Value = 1
$ws.Range("R14").Value = 0

# --- color ---
$ws.Range("U13").Value = "white"
$ws.Range("U14").Value = "light gray"

# --- date_complete ---
$ws.Range("AV13").Value = 44714
$ws.Range("AV14").Value = 44707

# --- gcms_peak ---
$ws.Range("AY13").Value = 7.630000114440918
$ws.Range("AY14").Value = 9.2399997711181641

# --- lab_num_substances_any ---
$ws.Range("BA13").Value = 7
$ws.Range("BA14").Value = 9

# --- lab_num_substances ---
$ws.Range("BB13").Value = 1
$ws.Range("BB14").Value = 5

# --- lab_fentanyl ---
$ws.Range("BC13").Value = 0
$ws.Range("BC14").Value = 1

# --- lab_fentanyl_any ---
$ws.Range("BD13").Value = 0
$ws.Range("BD14").Value = 1

# --- lab_xylazine ---
$ws.Range("BE13").Value = 0
$ws.Range("BE14").Value = 1

# --- lab_xylazine_any ---
$ws.Range("BF13").Value = 0
$ws.Range("BF14").Value = 1

# --- lab_cocaine ---
$ws.Range("BJ13").Value = 1
$ws.Range("BJ14").Value = 0

# --- lab_cocaine_any ---
$ws.Range("BK13").Value = 1
$ws.Range("BK14").Value = 0

# --- lab_levamisole_any ---
$ws.Range("BN13").Value = 1
$ws.Range("BN14").Value = 0

# --- lab_opioid ---
$ws.Range("BW13").Value = 0
$ws.Range("BW14").Value = 1

# --- lab_opioid_any ---
$ws.Range("BX13").Value = 0
$ws.Range("BX14").Value = 1

# --- lab_cocaine_impurities_any ---
$ws.Range("CI13").Value = 1
$ws.Range("CI14").Value = 0

# --- lab_heroin_impurities_any ---
$ws.Range("CK13").Value = 0
$ws.Range("CK14").Value = 1

# --- lab_fentanyl_impurities_any ---
$ws.Range("CM13").Value = 0
$ws.Range("CM14").Value = 1
